# Add a new "商品分類" (Product Category) column between the existing
# "檔名判斷" and "對應角度" columns, and record "運動鞋" (Sneakers) as the
# category value for the "_PR_" / "鞋盒" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B; the former column B ("對應角度" + its values)
# shifts right to become column C.
$ws.Columns("B:B").Insert()

# New header + value for the inserted column.
$ws.Range("B1").Value = "商品分類"
$ws.Range("B7").Value = "運動鞋"

# Match the fonts Excel applied to the new column's cells.
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 10

$ws.Range("B2").Font.Name = "Microsoft JhengHei"
$ws.Range("B2").Font.Color = 0

$ws.Range("B4").Font.Name = "Microsoft JhengHei"

$ws.Range("B7").Font.Name = "Microsoft JhengHei"
$ws.Range("B7").Font.Color = 0

# Keep the existing 15.75pt custom row height explicit on every row.
$ws.Rows("1:15").RowHeight = 15.75

# Leave the selection on the cell that was last edited.
$null = $ws.Range("B7").Select()
